# Add "Norway" and "Poland" market test-data sheets, cloned from the
# existing "Austria" sheet (same layout/column widths as the target
# sheets), positioned right after "Hungary".

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Austria")
$hungary  = $wb.Worksheets.Item("Hungary")

# --- Norway -------------------------------------------------------------
$template.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($hungary.Index + 1)
$norway.Name = "Norway"

$norway.Range("B4").Value = "NGC-2931/T3086/T3084"
$norway.Range("B2").Value = "Norway Market"

# Insert the extra "MZXSDR240" repeater row (between MZX64DR and MZXDR240)
$norway.Rows(15).Insert()
$norway.Range("A14").Copy()
$norway.Range("A15").PasteSpecial(-4122)
$norway.Range("A15").Value = "MZXSDR240"

# --- Poland ---------------------------------------------------------------
$template.Copy($null, $norway)
$poland = $wb.Worksheets.Item($norway.Index + 1)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/T3121/T3119"
$poland.Range("B2").Value = "Poland Market"

# Norway is the sheet that should end up active/selected (matches the
# authored workbook state), not Poland (the last-created sheet).
$norway.Activate()
